$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Vertical Offset" column (H) is being removed from the table.
# Columns to its right ("Points Density" in I, "End Time" in J) shift one
# column to the left (I -> H, J -> I), and the former last column (J) is
# cleared out entirely.
#
# Values are copied/shifted cell-by-cell (instead of using a whole-column
# Delete/shift) so that the per-column <col> width definitions stay tied to
# their original column index, matching how the sheet ended up after the
# edit.

$lastRow = 10
$colVerticalOffset = 8  # H
$colPointsDensity  = 9  # I
$colEndTime        = 10 # J

# Step 1: shift the "Points Density" and "End Time" values left by one
# column for every row (header row included).
for ($r = 1; $r -le $lastRow; $r++) {
    $pointsDensityValue = $ws.Cells.Item($r, $colPointsDensity).Value2
    $endTimeValue = $ws.Cells.Item($r, $colEndTime).Value2

    $ws.Cells.Item($r, $colVerticalOffset).Value2 = $pointsDensityValue
    $ws.Cells.Item($r, $colPointsDensity).Value2 = $endTimeValue
}

# Step 2: remove the now-empty trailing column's cells entirely (shifting
# cells left so nothing remains in column J).
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $colEndTime).Delete(-4159)  # xlShiftToLeft
}

# Reflect the selection the user ended up with after performing the edit.
$ws.Range("H1:I10").Select()

Write-Host "Removed 'Vertical Offset' column and shifted remaining columns left"
